$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (becomes the "Björksplintborre" record, previously in row 17) ---
$ws.Cells.Item(16, 1).Value = 130827881
$ws.Cells.Item(16, 2).Value = 8440
$ws.Cells.Item(16, 5).Value = 106554
$ws.Cells.Item(16, 6).Value = "Björksplintborre"
$ws.Cells.Item(16, 7).Value = "Scolytus ratzeburgii"
$ws.Cells.Item(16, 8).Value = "Janson, 1856"
$ws.Cells.Item(16, 17).Value = 344495
$ws.Cells.Item(16, 18).Value = 6433268
$ws.Cells.Item(16, 36).ClearContents()
$ws.Cells.Item(16, 37).ClearContents()
$ws.Cells.Item(16, 41).ClearContents()

# --- Row 17 (becomes the "Rostfläck" record, previously in row 16) ---
$ws.Cells.Item(17, 1).Value = 130827869
$ws.Cells.Item(17, 2).Value = 75222
$ws.Cells.Item(17, 5).Value = 6428
$ws.Cells.Item(17, 6).Value = "Rostfläck"
$ws.Cells.Item(17, 7).Value = "Arthonia vinosa"
$ws.Cells.Item(17, 8).Value = "Leight."
$ws.Cells.Item(17, 17).Value = 344518
$ws.Cells.Item(17, 18).Value = 6433262
$ws.Cells.Item(17, 36).Value = "gran"
$ws.Cells.Item(17, 37).Value = "Picea abies"
$ws.Cells.Item(17, 41).Value = "Picea abies"

# --- Row 19 (becomes the "Kungsfågel" record, previously in row 20) ---
$ws.Cells.Item(19, 1).Value = 130827872
$ws.Cells.Item(19, 2).Value = 58256
$ws.Cells.Item(19, 5).Value = 103015
$ws.Cells.Item(19, 6).Value = "Kungsfågel"
$ws.Cells.Item(19, 7).Value = "Regulus regulus"
$ws.Cells.Item(19, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(19, 17).Value = 344461
$ws.Cells.Item(19, 18).Value = 6433331
$ws.Cells.Item(19, 36).ClearContents()
$ws.Cells.Item(19, 37).ClearContents()
$ws.Cells.Item(19, 41).ClearContents()

# --- Row 20 (becomes the "Kattfotslav" record, previously in row 19) ---
$ws.Cells.Item(20, 1).Value = 130827874
$ws.Cells.Item(20, 2).Value = 75350
$ws.Cells.Item(20, 5).Value = 6426
$ws.Cells.Item(20, 6).Value = "Kattfotslav"
$ws.Cells.Item(20, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(20, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(20, 17).Value = 344473
$ws.Cells.Item(20, 18).Value = 6433342
$ws.Cells.Item(20, 36).Value = "gran"
$ws.Cells.Item(20, 37).Value = "Picea abies"
$ws.Cells.Item(20, 41).Value = "Picea abies"

# --- Row 21: only Taxonsorteringsordning (B) changes ---
$ws.Cells.Item(21, 2).Value = 97629

# --- Row 22: only Taxonsorteringsordning (B) changes ---
$ws.Cells.Item(22, 2).Value = 83207

# --- Row 23: only Taxonsorteringsordning (B) changes ---
$ws.Cells.Item(23, 2).Value = 75350

# --- Row 26 (becomes the "Kornig nållav" record, previously in row 27) ---
$ws.Cells.Item(26, 1).Value = 130827873
$ws.Cells.Item(26, 2).Value = 83209
$ws.Cells.Item(26, 5).Value = 306
$ws.Cells.Item(26, 6).Value = "Kornig nållav"
$ws.Cells.Item(26, 7).Value = "Chaenotheca chlorella"
$ws.Cells.Item(26, 8).Value = "(Ach.) Müll.Arg."
$ws.Cells.Item(26, 17).Value = 344451
$ws.Cells.Item(26, 18).Value = 6433334

# --- Row 27 (becomes the "Kattfotslav" record, previously in row 26) ---
$ws.Cells.Item(27, 1).Value = 130827876
$ws.Cells.Item(27, 2).Value = 75350
$ws.Cells.Item(27, 5).Value = 6426
$ws.Cells.Item(27, 6).Value = "Kattfotslav"
$ws.Cells.Item(27, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(27, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(27, 17).Value = 344449
$ws.Cells.Item(27, 18).Value = 6433318
